# Add two redundant TestCase fields ("test_case_predicate" and
# "test_case_input_id") to every sheet that models a TestCase (or a
# subclass of it). The new columns are inserted immediately after the
# existing "test_case_source" column and immediately before the "id"
# column, pushing id/name/description/tags (and anything else after
# them) two columns to the right.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "TestCase",
    "AcceptanceTestCase",
    "QuantitativeTestCase",
    "ComplianceTestCase",
    "KnowledgeGraphNavigationTestCase",
    "OneHopTestCase"
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the column whose header in row 1 is "test_case_source"; the
    # new columns are inserted right after it (i.e. right before "id").
    $lastCol = $ws.UsedRange.Columns.Count
    $insertCol = 0
    for ($c = 1; $c -le $lastCol; $c++) {
        $headerValue = $ws.Cells.Item(1, $c).Text
        if ($headerValue -eq "test_case_source") {
            $insertCol = $c + 1
            break
        }
    }

    # Insert two new (blank) columns at $insertCol, shifting the
    # existing id/name/description/tags columns two places right.
    $ws.Columns.Item($insertCol).Insert()
    $ws.Columns.Item($insertCol).Insert()

    $ws.Cells.Item(1, $insertCol).Value = "test_case_predicate"
    $ws.Cells.Item(1, $insertCol + 1).Value = "test_case_input_id"
}
